$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.339722156331853
$ws.Range("C2").Value = 0.7817454319255371

$ws.Range("B3").Value = 6.876393158390012
$ws.Range("C3").Value = 0.9017652232179459

$ws.Range("B4").Value = 4.22508622341286
$ws.Range("C4").Value = 0.7912978609182165

$ws.Range("B5").Value = 3.087154738082091
$ws.Range("C5").Value = 0.997972229015301

$ws.Range("B6").Value = 2.136267898757287
$ws.Range("C6").Value = 0.9767216901501417

$ws.Range("B7").Value = 1.772567729512532
$ws.Range("C7").Value = 0.999006176665434

$ws.Range("B8").Value = 2.21468910484377
$ws.Range("C8").Value = 0.9974631771063818

$ws.Range("B9").Value = 15.15625224194849
$ws.Range("C9").Value = 0.8188174775617363

$ws.Range("B10").Value = 1.586595234740817
$ws.Range("C10").Value = 0.9952704087474706
